$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels
$ws.Range("B1").Value = "applicants"
$ws.Range("C1").Value = "percentageAdmitted"
$ws.Range("D1").Value = "percEnrolled"
$ws.Range("E1").Value = "totalEnrolled"

# Fill column E with the computed "totalEnrolled" formula for each data row.
# E2 is entered as a standalone formula; E3:E14 are entered together as a
# range so they form a shared formula group (matching manual entry followed
# by a multi-cell fill/paste for the remaining rows).
$ws.Range("E2").Formula = "=INT(B2*C2*D2)"
$ws.Range("E3:E14").Formula = "=INT(B3*C3*D3)"

# Move the active selection to match the saved workbook state
$ws.Range("F13").Select() | Out-Null
